# Appends three new species-observation rows (62-64) to the "Artfynd"
# worksheet, matching rows already present in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62
$ws.Cells.Item(62, 1).Value = 131134870
$ws.Cells.Item(62, 2).Value = 56748
$ws.Cells.Item(62, 4).Value = "NT"
$ws.Cells.Item(62, 5).Value = 205998
$ws.Cells.Item(62, 6).Value = "Nordfladdermus"
$ws.Cells.Item(62, 7).Value = "Eptesicus nilssonii"
$ws.Cells.Item(62, 8).Value = "(A.Keyserling & Blasius, 1839)"
$ws.Cells.Item(62, 9).Value = ""
$ws.Cells.Item(62, 14).Value = "autobox"
$ws.Cells.Item(62, 16).Value = "Högåsen, Dlr"
$ws.Cells.Item(62, 17).Value = 557180
$ws.Cells.Item(62, 18).Value = 6710045
$ws.Cells.Item(62, 19).Value = 10
$ws.Cells.Item(62, 20).Value = "Dalarna"
$ws.Cells.Item(62, 21).Value = "Hedemora"
$ws.Cells.Item(62, 22).Value = "Dalarna"
$ws.Cells.Item(62, 23).Value = "Husby"
$ws.Cells.Item(62, 25).NumberFormat = "@"
$ws.Cells.Item(62, 25).Value = "2025-06-23"
$ws.Cells.Item(62, 27).NumberFormat = "@"
$ws.Cells.Item(62, 27).Value = "2025-07-31"
$ws.Cells.Item(62, 30).Value = $False
$ws.Cells.Item(62, 31).Value = $False
$ws.Cells.Item(62, 33).Value = $False
$ws.Cells.Item(62, 44).Value = ""
$ws.Cells.Item(62, 46).Value = ""
$ws.Cells.Item(62, 49).Value = "Emmy Ransgart"
$ws.Cells.Item(62, 50).Value = "Via Emmy Ransgart"
$ws.Cells.Item(62, 51).Value = ""

# Row 63
$ws.Cells.Item(63, 1).Value = 131134868
$ws.Cells.Item(63, 2).Value = 56755
$ws.Cells.Item(63, 4).Value = "LC"
$ws.Cells.Item(63, 5).Value = 205992
$ws.Cells.Item(63, 6).Value = "Vattenfladdermus"
$ws.Cells.Item(63, 7).Value = "Myotis daubentonii"
$ws.Cells.Item(63, 8).Value = "(Kuhl, 1817)"
$ws.Cells.Item(63, 9).Value = ""
$ws.Cells.Item(63, 14).Value = "autobox"
$ws.Cells.Item(63, 16).Value = "Högåsen, Dlr"
$ws.Cells.Item(63, 17).Value = 557180
$ws.Cells.Item(63, 18).Value = 6710045
$ws.Cells.Item(63, 19).Value = 10
$ws.Cells.Item(63, 20).Value = "Dalarna"
$ws.Cells.Item(63, 21).Value = "Hedemora"
$ws.Cells.Item(63, 22).Value = "Dalarna"
$ws.Cells.Item(63, 23).Value = "Husby"
$ws.Cells.Item(63, 25).NumberFormat = "@"
$ws.Cells.Item(63, 25).Value = "2025-06-23"
$ws.Cells.Item(63, 27).NumberFormat = "@"
$ws.Cells.Item(63, 27).Value = "2025-07-31"
$ws.Cells.Item(63, 30).Value = $False
$ws.Cells.Item(63, 31).Value = $False
$ws.Cells.Item(63, 33).Value = $False
$ws.Cells.Item(63, 44).Value = ""
$ws.Cells.Item(63, 46).Value = ""
$ws.Cells.Item(63, 49).Value = "Emmy Ransgart"
$ws.Cells.Item(63, 50).Value = "Via Emmy Ransgart"
$ws.Cells.Item(63, 51).Value = ""

# Row 64
$ws.Cells.Item(64, 1).Value = 131134867
$ws.Cells.Item(64, 2).Value = 56769
$ws.Cells.Item(64, 4).Value = "NT"
$ws.Cells.Item(64, 5).Value = 206002
$ws.Cells.Item(64, 6).Value = "Brunlångöra"
$ws.Cells.Item(64, 7).Value = "Plecotus auritus"
$ws.Cells.Item(64, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(64, 9).Value = ""
$ws.Cells.Item(64, 14).Value = "autobox"
$ws.Cells.Item(64, 16).Value = "Högåsen, Dlr"
$ws.Cells.Item(64, 17).Value = 557180
$ws.Cells.Item(64, 18).Value = 6710045
$ws.Cells.Item(64, 19).Value = 10
$ws.Cells.Item(64, 20).Value = "Dalarna"
$ws.Cells.Item(64, 21).Value = "Hedemora"
$ws.Cells.Item(64, 22).Value = "Dalarna"
$ws.Cells.Item(64, 23).Value = "Husby"
$ws.Cells.Item(64, 25).NumberFormat = "@"
$ws.Cells.Item(64, 25).Value = "2025-06-23"
$ws.Cells.Item(64, 27).NumberFormat = "@"
$ws.Cells.Item(64, 27).Value = "2025-07-31"
$ws.Cells.Item(64, 30).Value = $False
$ws.Cells.Item(64, 31).Value = $False
$ws.Cells.Item(64, 33).Value = $False
$ws.Cells.Item(64, 44).Value = ""
$ws.Cells.Item(64, 46).Value = ""
$ws.Cells.Item(64, 49).Value = "Emmy Ransgart"
$ws.Cells.Item(64, 50).Value = "Via Emmy Ransgart"
$ws.Cells.Item(64, 51).Value = ""
